$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.719288
$ws.Range("H2").Value = 11.157864
$ws.Range("I2").Value = 0.04235839908674209
$ws.Range("J2").Value = 0.04235839908674209
$ws.Range("O2").Value = 0.02773017886769741
$ws.Range("P2").Value = 0.02773017886769741
$ws.Range("Q2").Value = 0.1931847777706667
$ws.Range("R2").Value = 1.738662999936
$ws.Range("S2").Value = 0.001174605983224669
$ws.Range("T2").Value = 0.001174605983224669

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.719288
$ws.Range("H3").Value = 11.157864
$ws.Range("I3").Value = 0.04235839908674209
$ws.Range("J3").Value = 0.04235839908674209
$ws.Range("M3").Value = 1.821156333333333
$ws.Range("N3").Value = 5.463469
$ws.Range("O3").Value = 0.9722698211323025
$ws.Range("P3").Value = 0.9722698211323026
$ws.Range("Q3").Value = 6.773404896690667
$ws.Range("R3").Value = 60.960644070216
$ws.Range("S3").Value = 0.04118379310351742
$ws.Range("T3").Value = 0.04118379310351743

$ws.Range("I4").Value = 0.2979256989470644
$ws.Range("J4").Value = 0.2979256989470644
$ws.Range("O4").Value = 0.02773017886769741
$ws.Range("P4").Value = 0.02773017886769741
$ws.Range("S4").Value = 0.008261532921085866
$ws.Range("T4").Value = 0.008261532921085866

$ws.Range("I5").Value = 0.2979256989470644
$ws.Range("J5").Value = 0.2979256989470644
$ws.Range("M5").Value = 1.821156333333333
$ws.Range("N5").Value = 5.463469
$ws.Range("O5").Value = 0.9722698211323025
$ws.Range("P5").Value = 0.9722698211323026
$ws.Range("Q5").Value = 47.64040737152523
$ws.Range("R5").Value = 428.763666343727
$ws.Range("S5").Value = 0.2896641660259785
$ws.Range("T5").Value = 0.2896641660259786

$ws.Range("G6").Value = 6.299630666666666
$ws.Range("H6").Value = 18.898892
$ws.Range("I6").Value = 0.07174552491706633
$ws.Range("J6").Value = 0.07174552491706633
$ws.Range("O6").Value = 0.02773017886769741
$ws.Range("P6").Value = 0.02773017886769741
$ws.Range("Q6").Value = 0.3272112163342222
$ws.Range("R6").Value = 2.944900947008
$ws.Range("S6").Value = 0.001989516238907091
$ws.Range("T6").Value = 0.001989516238907091

$ws.Range("G7").Value = 6.299630666666666
$ws.Range("H7").Value = 18.898892
$ws.Range("I7").Value = 0.07174552491706633
$ws.Range("J7").Value = 0.07174552491706633
$ws.Range("M7").Value = 1.821156333333333
$ws.Range("N7").Value = 5.463469
$ws.Range("O7").Value = 0.9722698211323025
$ws.Range("P7").Value = 0.9722698211323026
$ws.Range("Q7").Value = 11.47261228626089
$ws.Range("R7").Value = 103.253510576348
$ws.Range("S7").Value = 0.06975600867815923
$ws.Range("T7").Value = 0.06975600867815925

$ws.Range("G8").Value = 51.62686066666667
$ws.Range("H8").Value = 154.880582
$ws.Range("I8").Value = 0.5879703770491272
$ws.Range("J8").Value = 0.5879703770491272
$ws.Range("O8").Value = 0.02773017886769741
$ws.Range("P8").Value = 0.02773017886769741
$ws.Range("Q8").Value = 2.681567978840889
$ws.Range("R8").Value = 24.134111809568
$ws.Range("S8").Value = 0.01630452372447979
$ws.Range("T8").Value = 0.01630452372447979

$ws.Range("G9").Value = 51.62686066666667
$ws.Range("H9").Value = 154.880582
$ws.Range("I9").Value = 0.5879703770491272
$ws.Range("J9").Value = 0.5879703770491272
$ws.Range("M9").Value = 1.821156333333333
$ws.Range("N9").Value = 5.463469
$ws.Range("O9").Value = 0.9722698211323025
$ws.Range("P9").Value = 0.9722698211323026
$ws.Range("Q9").Value = 94.02058427321755
$ws.Range("R9").Value = 846.1852584589581
$ws.Range("S9").Value = 0.5716658533246475
$ws.Range("T9").Value = 0.5716658533246475

Write-Output "applied changes"